$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cell values for rows 2-7 (sender/target cluster remap + new TPM-derived stats)
$ws.Range("D2").Value2 = "FAPs"
$ws.Range("G2").Value2 = 0.0007963333333333334
$ws.Range("H2").Value2 = 0.002389
$ws.Range("I2").Value2 = 0.0106576135689399
$ws.Range("J2").Value2 = 0.01065761356893991
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 6.072364333333334
$ws.Range("N2").Value2 = 18.217093
$ws.Range("O2").Value2 = 0.4407767221912973
$ws.Range("P2").Value2 = 0.4407767221912974
$ws.Range("Q2").Value2 = 0.004835626130777778
$ws.Range("R2").Value2 = 0.04352063517700001
$ws.Range("S2").Value2 = 0.004697627975298825
$ws.Range("T2").Value2 = 0.004697627975298826
$ws.Range("D3").Value2 = "MuSCs"
$ws.Range("G3").Value2 = 0.0007963333333333334
$ws.Range("H3").Value2 = 0.002389
$ws.Range("I3").Value2 = 0.0106576135689399
$ws.Range("J3").Value2 = 0.01065761356893991
$ws.Range("M3").Value2 = 7.704144333333335
$ws.Range("N3").Value2 = 23.112433
$ws.Range("O3").Value2 = 0.5592232778087027
$ws.Range("P3").Value2 = 0.5592232778087027
$ws.Range("Q3").Value2 = 0.006135066937444446
$ws.Range("R3").Value2 = 0.05521560243700001
$ws.Range("S3").Value2 = 0.00595998559364108
$ws.Range("T3").Value2 = 0.005959985593641081
$ws.Range("A4").Value2 = "FAPs"
$ws.Range("D4").Value2 = "FAPs"
$ws.Range("G4").Value2 = 0.05547833333333333
$ws.Range("H4").Value2 = 0.166435
$ws.Range("I4").Value2 = 0.7424863601283017
$ws.Range("J4").Value2 = 0.7424863601283018
$ws.Range("M4").Value2 = 6.072364333333334
$ws.Range("N4").Value2 = 18.217093
$ws.Range("O4").Value2 = 0.4407767221912973
$ws.Range("P4").Value2 = 0.4407767221912974
$ws.Range("Q4").Value2 = 0.3368846526061111
$ws.Range("R4").Value2 = 3.031961873455
$ws.Range("S4").Value2 = 0.3272707040891
$ws.Range("T4").Value2 = 0.3272707040891001
$ws.Range("D5").Value2 = "MuSCs"
$ws.Range("I5").Value2 = 0.7424863601283017
$ws.Range("J5").Value2 = 0.7424863601283018
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 7.704144333333335
$ws.Range("N5").Value2 = 23.112433
$ws.Range("O5").Value2 = 0.5592232778087027
$ws.Range("P5").Value2 = 0.5592232778087027
$ws.Range("Q5").Value2 = 0.4274130873727778
$ws.Range("R5").Value2 = 3.846717786355001
$ws.Range("S5").Value2 = 0.4152156560392018
$ws.Range("T5").Value2 = 0.4152156560392019
$ws.Range("A6").Value2 = "MuSCs"
$ws.Range("G6").Value2 = 0.018445
$ws.Range("H6").Value2 = 0.055335
$ws.Range("I6").Value2 = 0.2468560263027583
$ws.Range("J6").Value2 = 0.2468560263027583
$ws.Range("O6").Value2 = 0.4407767221912973
$ws.Range("P6").Value2 = 0.4407767221912974
$ws.Range("Q6").Value2 = 0.1120047601283333
$ws.Range("R6").Value2 = 1.008042841155
$ws.Range("S6").Value2 = 0.1088083901268985
$ws.Range("T6").Value2 = 0.1088083901268985
$ws.Range("A7").Value2 = "MuSCs"
$ws.Range("G7").Value2 = 0.018445
$ws.Range("H7").Value2 = 0.055335
$ws.Range("I7").Value2 = 0.2468560263027583
$ws.Range("J7").Value2 = 0.2468560263027583
$ws.Range("M7").Value2 = 7.704144333333335
$ws.Range("N7").Value2 = 23.112433
$ws.Range("O7").Value2 = 0.5592232778087027
$ws.Range("P7").Value2 = 0.5592232778087027
$ws.Range("Q7").Value2 = 0.1421029422283334
$ws.Range("R7").Value2 = 1.278926480055
$ws.Range("S7").Value2 = 0.1380476361758598
$ws.Range("T7").Value2 = 0.1380476361758599

# Remove the now-obsolete rows (old rows 8-10, MuSCs-sourced edges to ECs / shifted pairs)
$ws.Rows("8:10").Delete() | Out-Null

Write-Output "Colq-Musk sheet updated with new TPM values"